# Aggiornamento dati fino al 09/09 (incluso)
# Appends rows 367..374 (dates 44441..44448, i.e. 2021-09-02 .. 2021-09-09)
# to the existing daily-data table on the active sheet, matching the
# style/format already used for the preceding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 366
$newLastRow = 374

# Copy formatting (number format, font, borders, alignment) from the last
# existing data row down onto the new rows, so new cells look the same as
# the rest of the table (e.g. the date style in column A).
$ws.Range("A" + $lastRow + ":D" + $lastRow).Copy()
$ws.Range("A" + ($lastRow + 1) + ":D" + $newLastRow).PasteSpecial(-4122)
$excel.CutCopyMode = 0

$dates = @(44441, 44442, 44443, 44444, 44445, 44446, 44447, 44448)
$nuoviPos = @(0, 0, 0, 0, 0, 0, 0, 2)
$sommaMobile = @(1, 1, 0, 0, 0, 0, 0, 2)
$sommaMobile100k = @(16.02307322544464, 16.02307322544464, 0, 0, 0, 0, 0, 32.04614645088928)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $lastRow + 1 + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $nuoviPos[$i]
    $ws.Cells.Item($r, 3).Value = $sommaMobile[$i]
    $ws.Cells.Item($r, 4).Value = $sommaMobile100k[$i]
}
